$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Batch registration mode: the "postal code" column (I) now holds a real
# numeric postal code instead of reusing the shared "address" string.
$ws.Range("I2").Value = 12345

# Update the active selection on the sheet to I3
$ws.Range("I3").Select()
